$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Sort data rows A2:M9 by column B (RossED50) ascending
$rng = $ws.Range("A2:M9")
$key = $ws.Range("B2:B9")
$rng.Sort($key, 1)

# 2. Insert a new blank column at L (12); this shifts the old L (CKRank) -> M
#    and the old M (USRank) -> N, leaving a blank column L for the new RossCKRank data.
$ws.Columns.Item(12).Insert()

# 3. New header cells. Insert USvsRossRank (O1) before RossCKRank (L1) so that the
#    shared-string table ends up with the same ordering as the target workbook.
$ws.Range("O1").Value2 = "USvsRossRank"
$ws.Range("L1").Value2 = "RossCKRank"
$ws.Range("P1").Value2 = "RossvsCKRank"

# 4. Fill new L column (RossCKRank) for the 5 genotypes that have a CK value
$ws.Range("L4").Value2 = 1
$ws.Range("L5").Value2 = 2
$ws.Range("L6").Value2 = 3
$ws.Range("L7").Value2 = 4
$ws.Range("L8").Value2 = 5

# 5. New O column (USvsRossRank) = ABS(USRank - RossRank) for every data row
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 15).Formula = "=ABS(N$r-K$r)"
}

# 6. New P column (RossvsCKRank) = ABS(RossCKRank - CKRank) for rows 4-8, as a shared formula
$ws.Range("P4").Formula = "=ABS(L4-M4)"
$ws.Range("P5:P8").Formula = "=ABS(L5-M5)"

# 7. New row 10 summary formulas
$ws.Range("B10").Formula = "=MAX(B2:B9)-MIN(B2:B9)"
$ws.Range("F10").Formula = "=MAX(F2:F9)-MIN(F2:F9)"
$ws.Range("G10").Formula = "=MAX(G2:G9)-MIN(G2:G9)"
$ws.Range("O10").Formula = "=AVERAGE(O2:O9)"
$ws.Range("P10").Formula = "=AVERAGE(P4:P8)"

# 8. Column width for the new O column (USvsRossRank)
$ws.Columns.Item(15).ColumnWidth = 11.6640625

# 9. Re-apply the sort through the Sort object so that the persisted sortState
#    element reflects the new range/condition (A2:O9 sorted by B2:B9)
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B9")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:O9"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# 10. Selection / active cell, matching the saved view in the target workbook
$ws.Range("M5").Select()

Write-Host "done"
